# teacher services done correctly
# Updates the "Teacher Details " sheet in the TeacherDB workbook:
#  - fixes teacher1's email/age/phone on row 2
#  - replaces the duplicated row 3 with teacher2's corrected data
#  - replaces teacher2's old row 4 with a new teacher10 entry
#  - appends a brand new teacher5 entry on row 5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: teacher1 -- fix age, phone and email ---
$ws.Range("C2").Value = 30
$ws.Range("E2").Value = 1111111111
$ws.Range("F2").Value = "teacher1@gmail.com"

# --- Row 3: was a duplicate of teacher1, now holds teacher2's data ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "teacher2"
$ws.Range("D3").Value = "5->Math"
$ws.Range("E3").Value = 2222222222
$ws.Range("F3").Value = "teacher2@gmail.com"
$ws.Range("G3").Value = "2 cairo street"

# --- Row 4: was teacher2, now holds the new teacher10 entry ---
$ws.Range("A4").Value = 15
$ws.Range("B4").Value = "teacher10"
$ws.Range("C4").Value = 40
$ws.Range("D4").Value = "6->Arabic"
$ws.Range("E4").Value = 1010101010
$ws.Range("F4").Value = "teacher10@"
$ws.Range("G4").Value = " cairo street"

# --- Row 5: brand new teacher5 entry ---
$ws.Range("A5").Value = 20
$ws.Range("D5").Value = "1->histroy"
$ws.Range("F5").Value = "teacher5@gmail.com"
$ws.Range("G5").Value = "street "
$ws.Range("B5").Value = "teacher5"
$ws.Range("C5").Value = 33
$ws.Range("E5").Value = 2020202020
$ws.Range("H5").Value = $false

# Select the newly added row, same as the author leaving the cursor there
$ws.Rows(5).Select()
